# Update column G ("K" = Strike count) values for the 2023 sulser_cole save data.
# Per commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" - here we only need to rewrite the recalculated K values
# for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 0
    7  = 1
    8  = 2
    9  = 3
    10 = 1
    11 = 0
    13 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
